$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the File Location text for "Subscription Model" row (was under Managers/Components, now just Components)
$ws.Range("C5").Value = " Components/Subscription.cpp"

# Fix typo: "Collision resonse" -> "Collision resolution"
$ws.Range("A14").Value = "Collision resolution"

# Update last active selected cell to A12
$ws.Range("A12").Select()

$wb.Save()
